$wb = $excel.ActiveWorkbook
# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H21").Value = 32679.666
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 14010
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 14010
$ws.Range("M21").Value = -69551
$ws.Range("N21").Value = -14946
$ws.Range("H23").Value = 32679.666
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 14010
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 14010
$ws.Range("M23").Value = -69785
$ws.Range("N23").Value = -14478
$ws.Range("H32").Value = 684.6
$ws.Range("J32").Value = 713.3333
$ws.Range("L32").Value = 713.3333
$ws.Range("N32").Value = -1365.3333
$ws.Range("H33").Value = 238.55556
$ws.Range("I33").Value = 259.6
$ws.Range("J33").Value = 212.25
$ws.Range("K33").Value = 259.6
$ws.Range("L33").Value = 212.25
$ws.Range("M33").Value = -30.60000000000002
$ws.Range("N33").Value = -670.25
$ws.Range("H39").Value = 900
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 900
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 2700
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -3292
$ws.Range("H51").Value = 3128.6667
$ws.Range("I51").Value = 3350
$ws.Range("J51").Value = 3105.3684
$ws.Range("K51").Value = 3350
$ws.Range("L51").Value = 3105.3684
$ws.Range("M51").Value = -2866
$ws.Range("N51").Value = -4073.3684
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H64").Value = 166671660
$ws.Range("I64").Value = 1000000000
$ws.Range("J64").Value = 5992
$ws.Range("K64").Value = 1000000000
$ws.Range("L64").Value = 5992
$ws.Range("M64").Value = -999999752
$ws.Range("N64").Value = -6488
$ws.Range("H67").Value = 166671660
$ws.Range("I67").Value = 1000000000
$ws.Range("J67").Value = 5992
$ws.Range("K67").Value = 1000000000
$ws.Range("L67").Value = 5992
$ws.Range("M67").Value = -999999142
$ws.Range("N67").Value = -7708
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = $null
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = $null
$ws.Range("H96").Value = 512
$ws.Range("I96").Value = 464.25
$ws.Range("J96").Value = 766.6667
$ws.Range("K96").Value = 1392.75
$ws.Range("L96").Value = 2300.0001
$ws.Range("M96").Value = -19.75
$ws.Range("N96").Value = -5046.0001
$ws.Range("H100").Value = 126972.5
$ws.Range("I100").Value = 167636.67
$ws.Range("J100").Value = 4980
$ws.Range("K100").Value = 167636.67
$ws.Range("L100").Value = 4980
$ws.Range("M100").Value = -167095.67
$ws.Range("N100").Value = -6062
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null
$ws.Range("H103").Value = 377
$ws.Range("I103").Value = 304
$ws.Range("J103").Value = 450
$ws.Range("K103").Value = 912
$ws.Range("L103").Value = 1350
$ws.Range("M103").Value = -326
$ws.Range("N103").Value = -2522
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
$ws.Range("H107").Value = 292.85715
$ws.Range("I107").Value = 292.85715
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 292.85715
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1627.14285
$ws.Range("N107").Value = $null
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H109").Value = 48000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 48000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 48000
$ws.Range("M109").Value = $null
$ws.Range("N109").Value = -50774
$ws.Range("H132").Value = 2571.2466
$ws.Range("I132").Value = 2254.547
$ws.Range("J132").Value = 3410.5
$ws.Range("K132").Value = 6763.641
$ws.Range("L132").Value = 10231.5
$ws.Range("M132").Value = -4233.641
$ws.Range("N132").Value = -15291.5
$ws.Range("H135").Value = 1122.0322
$ws.Range("I135").Value = 468
$ws.Range("J135").Value = 1735.1875
$ws.Range("K135").Value = 4212
$ws.Range("L135").Value = 15616.6875
$ws.Range("M135").Value = -1677
$ws.Range("N135").Value = -20686.6875
$ws.Range("H138").Value = 7996.8335
$ws.Range("I138").Value = 2264.3076
$ws.Range("J138").Value = 14771.637
$ws.Range("K138").Value = 6792.9228
$ws.Range("L138").Value = 44314.911
$ws.Range("M138").Value = -1652.9228
$ws.Range("N138").Value = -54594.911
$ws.Range("H141").Value = 9685.214
$ws.Range("I141").Value = 1853
$ws.Range("J141").Value = 38403.332
$ws.Range("K141").Value = 5559
$ws.Range("L141").Value = 115209.996
$ws.Range("M141").Value = -379
$ws.Range("N141").Value = -125569.996

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 49.973682
$ws.Range("I7").Value = 25
$ws.Range("J7").Value = 74.947365
$ws.Range("K7").Value = 25
$ws.Range("L7").Value = 74.947365
$ws.Range("M7").Value = 88
$ws.Range("N7").Value = -300.947365

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 6847.8
$ws.Range("I110").Value = 4800
$ws.Range("J110").Value = 15039
$ws.Range("K110").Value = 14400
$ws.Range("L110").Value = 45117
$ws.Range("M110").Value = -10310
$ws.Range("N110").Value = -53297

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6022.1763
$ws.Range("I132").Value = 6783.4165
$ws.Range("J132").Value = 4195.2
$ws.Range("K132").Value = 20350.2495
$ws.Range("L132").Value = 12585.6
$ws.Range("M132").Value = -17645.6
